$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (M2:T2)
$ws.Range("M2").Value = 2.425633666666667
$ws.Range("N2").Value = 7.276901000000001
$ws.Range("O2").Value = 0.0662600404061536
$ws.Range("P2").Value = 0.06626004040615362
$ws.Range("Q2").Value = 0.9610586979808889
$ws.Range("R2").Value = 8.649528281828001
$ws.Range("S2").Value = 0.0662600404061536
$ws.Range("T2").Value = 0.06626004040615362

# Row 3 updates (O3, P3, S3, T3)
$ws.Range("O3").Value = 0.4234968256437875
$ws.Range("P3").Value = 0.4234968256437876
$ws.Range("S3").Value = 0.4234968256437875
$ws.Range("T3").Value = 0.4234968256437876

# Row 4 updates (M4:T4)
$ws.Range("M4").Value = 18.67887366666667
$ws.Range("N4").Value = 56.03662100000001
$ws.Range("O4").Value = 0.5102431339500588
$ws.Range("P4").Value = 0.5102431339500588
$ws.Range("Q4").Value = 7.400744082887558
$ws.Range("R4").Value = 66.60669674598802
$ws.Range("S4").Value = 0.5102431339500588
$ws.Range("T4").Value = 0.5102431339500588
